# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet right after "总计" (i.e. before the
# existing "2022-Q2" sheet), fills it with its fund-holding data, and
# records the new quarter's summary row (count=2, value=0.02) at the top
# of the "总计" sheet's data table, shifting the other rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" sheet in the right place.
#    NOTE: the worksheet reference used as the insertion anchor becomes
#    stale immediately after Add(), so every sheet needed afterwards is
#    (re)fetched fresh, by position, once the sheet list is final.
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item(2)          # currently "2022-Q2"
$q3 = $wb.Worksheets.Add($anchor)         # new sheet placed right before it
$q3.Name = "2022-Q3"

$total = $wb.Worksheets.Item(1)           # "总计"

# ---------------------------------------------------------------------
# 2) Update the "总计" (total) sheet: shift existing data rows down one
#    row and write the new 2022-Q3 summary as the new row 2.
# ---------------------------------------------------------------------
$total.Range("A2:D4").Copy($total.Range("A3:D5"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.02

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# ---------------------------------------------------------------------
# 3) Populate the new "2022-Q3" sheet with its fund-holding table.
#    Columns B-G hold text values (number-look-alikes kept as text,
#    matching the other quarter sheets), column H holds real numbers.
# ---------------------------------------------------------------------
$headerStyleSrc = $total.Range("B1:D1")   # bold/bordered header style
$headerStyleSrc.Copy($q3.Range("B1:H1"))
$aColStyleSrc = $total.Range("A2:A3")
$aColStyleSrc.Copy($q3.Range("A2:A3"))

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3.Range("B2:G3").NumberFormat = "@"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "009999"
$q3.Range("C2").Value = "东方中国红利混合"
$q3.Range("D2").Value = "0.51"
$q3.Range("E2").Value = "79.37"
$q3.Range("F2").Value = "4.31"
$q3.Range("G2").Value = "0.0220"
$q3.Range("H2").Value = 4

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "001849"
$q3.Range("C3").Value = "前海开源强势共识100强等权重股票"
$q3.Range("D3").Value = "0.11"
$q3.Range("E3").Value = "91.68"
$q3.Range("F3").Value = "1.00"
$q3.Range("G3").Value = "0.0011"
$q3.Range("H3").Value = 10
